$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force number-formatted text cells to remain text (avoid Excel auto-converting
# strings like "1.00" or "10.00" into numeric values that lose formatting).
$textCells = @("D5","D9","D10","D12","D14","D16","D19","D20","D22","D25","D26","D28","D29","D31","D35","D38","D39","D40","D41","D42","D43","D46","D47","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data (price & volume refresh, plus two row swaps)
$ws.Range("D2").Value = "63.642.36"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "3.484.53"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "583.27"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("D7").Value = "3.485.35"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.485"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "4.068.88"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "27.41"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.455.57"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "63.745.98"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "10.00"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "14.32"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "382.84"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "3.623.24"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "72.98"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").Value = "1.58"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").Value = "3.492.71"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "23.48"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  +3.09%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.57"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Value = "161.15"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").Value = "0.0798"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Value = "26.58"
$ws.Range("E42").Value = "  +6.43%  "
$ws.Range("D43").Value = "0.809"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "41.45"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "6.83"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.414.14"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "0.897"
$ws.Range("E51").Value = "  +1.30%  "
